$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Cell values scraped from coinranking.com; Price (D) and Volume(1h) (E)
# columns are stored as plain text in the source sheet, so assign them
# as strings. A few Price values round-trip through Excel's automatic
# number detection and lose a trailing zero (e.g. '3.00' -> 3); for
# those we use a leading apostrophe to force literal text, exactly as
# typing them into the grid would.

$ws.Range('D2').Value = '44.343.65'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.258.20'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '316.65'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').Value = '99.99'
$ws.Range('E6').Value = '  -4.77%  '
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  -1.99%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.539'
$ws.Range('E9').Value = '  -5.53%  '
$ws.Range('D10').Value = '36.26'
$ws.Range('E10').Value = '  -6.17%  '
$ws.Range('D11').Value = '0.0827'
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('D12').Value = '7.42'
$ws.Range('E12').Value = '  -5.68%  '
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('D14').Value = '2.598.13'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').Value = '''0.850'
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = '2.242.80'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '14.06'
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('D18').Value = '44.078.10'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('D19').Value = '13.22'
$ws.Range('E19').Value = '  -5.35%  '
$ws.Range('D20').Value = '0.0₃0986'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').Value = '6.38'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '65.68'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').Value = '240.49'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').Value = '''3.00'
$ws.Range('E24').Value = '  -6.24%  '
$ws.Range('D25').Value = '2.05'
$ws.Range('E25').Value = '  -7.78%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').Value = '10.18'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = '38.07'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('E29').Value = '  -3.79%  '
$ws.Range('D30').Value = '6.06'
$ws.Range('E30').Value = '  -6.65%  '
$ws.Range('D31').Value = '''20.20'
$ws.Range('E31').Value = '  -2.23%  '
$ws.Range('D32').Value = '156.28'
$ws.Range('E32').Value = '  -3.40%  '
$ws.Range('D33').Value = '0.0844'
$ws.Range('E33').Value = '  -4.52%  '
$ws.Range('D34').Value = '3.46'
$ws.Range('E34').Value = '  +10.38%  '
$ws.Range('E35').Value = '  -3.36%  '
$ws.Range('D36').Value = '0.112'
$ws.Range('E36').Value = '  -3.98%  '
$ws.Range('D37').Value = '1.92'
$ws.Range('E37').Value = '  -4.64%  '
$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').Value = '15.33'
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('D40').Value = '3.55'
$ws.Range('E40').Value = '  -9.85%  '
$ws.Range('E41').Value = '  -11.56%  '
$ws.Range('D42').Value = '''0.0310'
$ws.Range('E42').Value = '  -5.67%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '1.725.94'
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('D45').Value = '84.41'
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').Value = '0.197'
$ws.Range('E46').Value = '  -5.73%  '
$ws.Range('D47').Value = '5.22'
$ws.Range('E47').Value = '  -4.27%  '
$ws.Range('D48').Value = '102.47'
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '56.89'
$ws.Range('E49').Value = '  -6.13%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.63'
$ws.Range('E50').Value = '  -4.32%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = '70.79'
$ws.Range('E51').Value = '  -5.69%  '
